$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B, C, D, E (G is the sum of these four)
$newB = 3.272327238179451
$newC = 1.626987699542094
$newE = 0.5333859586016987

# D differs for row 4 vs the other rows
$newD_default = 0.1496068669990043
$newD_row4 = 0.7210945179870265

# Corresponding sums (G column), taken directly from target data
$newG_default = 5.582307763322248
$newG_row4 = 6.15379541431027

$rows = 2,3,4,5
foreach ($r in $rows) {
    $ws.Cells.Item($r, 2).Value = $newB
    $ws.Cells.Item($r, 3).Value = $newC
    if ($r -eq 4) {
        $ws.Cells.Item($r, 4).Value = $newD_row4
        $ws.Cells.Item($r, 7).Value = $newG_row4
    } else {
        $ws.Cells.Item($r, 4).Value = $newD_default
        $ws.Cells.Item($r, 7).Value = $newG_default
    }
    $ws.Cells.Item($r, 5).Value = $newE
}
